# fix(publipostage): Correct status name
#
# "bleu" -> "noir" (statut_label) and the related statut_name wording
# ("... posté" -> "... postés ou publiés") are corrected throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use whole-cell exact matches so the longer "... dans les 36/12 mois"
# variants (which share a common prefix) are not clobbered by the shorter
# replacement strings.
$xlWhole = 2

$ws.Cells.Replace("bleu", "noir", $xlWhole)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $xlWhole)
